# Scheduled-runner refresh of market-price-derived Leve profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all crafting-class
# sheets. Values below are the refreshed snapshot; two cells (ARM!M43 and
# GSM!N102) no longer have a computed profit and are cleared instead.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9262557
$ws.Range("I51").Value = 2962.3333
$ws.Range("J51").Value = 13892354
$ws.Range("K51").Value = 2962.3333
$ws.Range("L51").Value = 13892354
$ws.Range("M51").Value = -2478.3333
$ws.Range("N51").Value = -13893322
$ws.Range("H74").Value = 7381.2856
$ws.Range("I74").Value = 3888.3333
$ws.Range("J74").Value = 10001
$ws.Range("K74").Value = 3888.3333
$ws.Range("L74").Value = 10001
$ws.Range("M74").Value = -2952.3333
$ws.Range("N74").Value = -11873
$ws.Range("H77").Value = 7381.2856
$ws.Range("I77").Value = 3888.3333
$ws.Range("J77").Value = 10001
$ws.Range("K77").Value = 19441.6665
$ws.Range("L77").Value = 50005
$ws.Range("M77").Value = -14761.6665
$ws.Range("N77").Value = -59365
$ws.Range("H86").Value = 2390.75
$ws.Range("J86").Value = 563
$ws.Range("L86").Value = 563
$ws.Range("N86").Value = -2809
$ws.Range("H89").Value = 2390.75
$ws.Range("J89").Value = 563
$ws.Range("L89").Value = 2815
$ws.Range("N89").Value = -14047
$ws.Range("H99").Value = 166667820
$ws.Range("I99").Value = 1416.6666
$ws.Range("J99").Value = 333334240
$ws.Range("K99").Value = 4249.9998
$ws.Range("L99").Value = 1000002720
$ws.Range("M99").Value = -2751.9998
$ws.Range("N99").Value = -1000005716
$ws.Range("H112").Value = 2018.4348
$ws.Range("J112").Value = 2097.1052
$ws.Range("L112").Value = 6291.3156
$ws.Range("N112").Value = -8507.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 119687.18
$ws.Range("I4").Value = 183547.45
$ws.Range("J4").Value = 2610
$ws.Range("K4").Value = 183547.45
$ws.Range("L4").Value = 2610
$ws.Range("M4").Value = -183431.45
$ws.Range("N4").Value = -2842
$ws.Range("H32").Value = 10345.296
$ws.Range("I32").Value = 8203.3125
$ws.Range("K32").Value = 8203.3125
$ws.Range("M32").Value = -7916.3125
$ws.Range("H43").Value = 10365.6
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10365.6
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 10365.6
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -10991.6
$ws.Range("H63").Value = 7545.364
$ws.Range("I63").Value = 4749.75
$ws.Range("J63").Value = 9142.857
$ws.Range("K63").Value = 4749.75
$ws.Range("L63").Value = 9142.857
$ws.Range("M63").Value = -4063.75
$ws.Range("N63").Value = -10514.857
$ws.Range("H66").Value = 7545.364
$ws.Range("I66").Value = 4749.75
$ws.Range("J66").Value = 9142.857
$ws.Range("K66").Value = 23748.75
$ws.Range("L66").Value = 45714.285
$ws.Range("M66").Value = -20316.75
$ws.Range("N66").Value = -52578.285
$ws.Range("H110").Value = 3489.9167
$ws.Range("I110").Value = 2292.95
$ws.Range("K110").Value = 2292.95
$ws.Range("M110").Value = -247.9499999999998
$ws.Range("H122").Value = 2145.2222
$ws.Range("I122").Value = 2180.7307
$ws.Range("K122").Value = 6542.1921
$ws.Range("M122").Value = -4092.1921
$ws.Range("H132").Value = 7184.186
$ws.Range("I132").Value = 7701.2163
$ws.Range("K132").Value = 23103.6489
$ws.Range("M132").Value = -20573.6489
$ws.Range("H133").Value = 65332.168
$ws.Range("J133").Value = 65332.168
$ws.Range("L133").Value = 65332.168
$ws.Range("N133").Value = -70392.16800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17863716
$ws.Range("I20").Value = 23817652
$ws.Range("K20").Value = 23817652
$ws.Range("M20").Value = -23817405
$ws.Range("H82").Value = 20632.455
$ws.Range("J82").Value = 23994.75
$ws.Range("L82").Value = 23994.75
$ws.Range("N82").Value = -24760.75
$ws.Range("H85").Value = 20632.455
$ws.Range("J85").Value = 23994.75
$ws.Range("L85").Value = 23994.75
$ws.Range("N85").Value = -26646.75
$ws.Range("H99").Value = 2458.111
$ws.Range("I99").Value = 2319.2917
$ws.Range("J99").Value = 3568.6667
$ws.Range("K99").Value = 2319.2917
$ws.Range("L99").Value = 3568.6667
$ws.Range("M99").Value = -821.2917000000002
$ws.Range("N99").Value = -6564.6667
$ws.Range("H134").Value = 85541.336
$ws.Range("I134").Value = 99072.19500000001
$ws.Range("K134").Value = 297216.585
$ws.Range("M134").Value = -294681.585

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1858.3334
$ws.Range("I16").Value = 1649.875
$ws.Range("J16").Value = 2096.5715
$ws.Range("K16").Value = 1649.875
$ws.Range("L16").Value = 2096.5715
$ws.Range("M16").Value = -1362.875
$ws.Range("N16").Value = -2670.5715
$ws.Range("H31").Value = 4323.383
$ws.Range("I31").Value = 3011.9583
$ws.Range("J31").Value = 5691.826
$ws.Range("K31").Value = 3011.9583
$ws.Range("L31").Value = 5691.826
$ws.Range("M31").Value = -2716.9583
$ws.Range("N31").Value = -6281.826
$ws.Range("H34").Value = 4323.383
$ws.Range("I34").Value = 3011.9583
$ws.Range("J34").Value = 5691.826
$ws.Range("K34").Value = 3011.9583
$ws.Range("L34").Value = 5691.826
$ws.Range("M34").Value = -2809.9583
$ws.Range("N34").Value = -6095.826
$ws.Range("H68").Value = 21979.9
$ws.Range("J68").Value = 21979.9
$ws.Range("L68").Value = 21979.9
$ws.Range("N68").Value = -23477.9
$ws.Range("H71").Value = 21979.9
$ws.Range("J71").Value = 21979.9
$ws.Range("L71").Value = 65939.70000000001
$ws.Range("N71").Value = -73427.70000000001
$ws.Range("H74").Value = 58307.5
$ws.Range("J74").Value = 58119.445
$ws.Range("L74").Value = 58119.445
$ws.Range("N74").Value = -59867.445
$ws.Range("H77").Value = 58307.5
$ws.Range("J77").Value = 58119.445
$ws.Range("L77").Value = 174358.335
$ws.Range("N77").Value = -183094.335
$ws.Range("H113").Value = 1858.3334
$ws.Range("I113").Value = 1649.875
$ws.Range("J113").Value = 2096.5715
$ws.Range("K113").Value = 1649.875
$ws.Range("L113").Value = 2096.5715
$ws.Range("M113").Value = 520.125
$ws.Range("N113").Value = -6436.5715
$ws.Range("H132").Value = 2562.0571
$ws.Range("I132").Value = 2107.3447
$ws.Range("K132").Value = 6322.034100000001
$ws.Range("M132").Value = -3792.034100000001
$ws.Range("H141").Value = 388639.22
$ws.Range("J141").Value = 388639.22
$ws.Range("L141").Value = 388639.22
$ws.Range("N141").Value = -398999.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 166699580
$ws.Range("J37").Value = 166699580
$ws.Range("L37").Value = 500098740
$ws.Range("N37").Value = -500098964
$ws.Range("H117").Value = 15880181
$ws.Range("I117").Value = 4692.8
$ws.Range("J117").Value = 20841272
$ws.Range("K117").Value = 14078.4
$ws.Range("L117").Value = 62523816
$ws.Range("M117").Value = -10636.4
$ws.Range("N117").Value = -62530700
$ws.Range("H122").Value = 1058.7142
$ws.Range("I122").Value = 882.4
$ws.Range("J122").Value = 1499.5
$ws.Range("K122").Value = 7941.599999999999
$ws.Range("L122").Value = 13495.5
$ws.Range("M122").Value = -5491.599999999999
$ws.Range("N122").Value = -18395.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3024.3547
$ws.Range("I102").Value = 3024.3547
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3024.3547
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1402.3547
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 3982.5833
$ws.Range("J122").Value = 3677.5
$ws.Range("L122").Value = 11032.5
$ws.Range("N122").Value = -15932.5
$ws.Range("H126").Value = 2320.1333
$ws.Range("I126").Value = 2172.9092
$ws.Range("J126").Value = 2725
$ws.Range("K126").Value = 6518.7276
$ws.Range("L126").Value = 8175
$ws.Range("M126").Value = -4048.7276
$ws.Range("N126").Value = -13115

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12273.538
$ws.Range("I40").Value = 14730.667
$ws.Range("K40").Value = 14730.667
$ws.Range("M40").Value = -14594.667
$ws.Range("H43").Value = 251734.62
$ws.Range("J43").Value = 251734.62
$ws.Range("L43").Value = 251734.62
$ws.Range("N43").Value = -252120.62
$ws.Range("H46").Value = 22412.182
$ws.Range("I46").Value = 34928.54
$ws.Range("K46").Value = 34928.54
$ws.Range("M46").Value = -34740.54
$ws.Range("H132").Value = 6437.385
$ws.Range("I132").Value = 6521.4443
$ws.Range("J132").Value = 6248.25
$ws.Range("K132").Value = 19564.3329
$ws.Range("L132").Value = 18744.75
$ws.Range("M132").Value = -17034.3329
$ws.Range("N132").Value = -23804.75
$ws.Range("H133").Value = 75867.336
$ws.Range("J133").Value = 75867.336
$ws.Range("L133").Value = 75867.336
$ws.Range("N133").Value = -80927.336
$ws.Range("H134").Value = 56742
$ws.Range("J134").Value = 56742
$ws.Range("L134").Value = 56742
$ws.Range("N134").Value = -66882

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 21856.143
$ws.Range("J4").Value = 22166.334
$ws.Range("L4").Value = 22166.334
$ws.Range("N4").Value = -22392.334
$ws.Range("H81").Value = 2301.6667
$ws.Range("I81").Value = 952.5
$ws.Range("K81").Value = 1905
$ws.Range("M81").Value = -844
$ws.Range("H84").Value = 2301.6667
$ws.Range("I84").Value = 952.5
$ws.Range("K84").Value = 9525
$ws.Range("M84").Value = -4221
$ws.Range("H96").Value = 2716.4092
$ws.Range("I96").Value = 1730.5714
$ws.Range("K96").Value = 1730.5714
$ws.Range("M96").Value = -357.5714
$ws.Range("H122").Value = 2372.3462
$ws.Range("I122").Value = 2403.8096
$ws.Range("K122").Value = 7211.4288
$ws.Range("M122").Value = -4761.4288
$ws.Range("H132").Value = 695219.25
$ws.Range("I132").Value = 964793.9399999999
$ws.Range("J132").Value = 4434.1875
$ws.Range("K132").Value = 2894381.82
$ws.Range("L132").Value = 13302.5625
$ws.Range("M132").Value = -2891851.82
$ws.Range("N132").Value = -18362.5625
